# The deck currently applies the "Integral" design (Red Violet color
# variant) as its one-and-only slide master/theme. This edit swaps the
# applied color variant back to the built-in "Office" palette - i.e. the
# Design > Variants > Colors gallery pick goes from "Red Violet" to
# "Office", while the font scheme / effect scheme (and everything else
# about the design) stay exactly as they were.
#
# All slides share the single slide master, so rewriting the theme's
# color scheme through any one slide's ThemeColorScheme updates the
# shared theme part for the whole deck.

$p = $ppt.ActivePresentation

function ConvertTo-ComRGB([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# Office theme color scheme, in the DrawingML clrScheme slot order:
# dk1, lt1, dk2, lt2, accent1, accent2, accent3, accent4, accent5,
# accent6, hlink, folHlink.
$officeColors = @(
    "000000", # dk1
    "FFFFFF", # lt1
    "44546A", # dk2
    "E7E6E6", # lt2
    "5B9BD5", # accent1
    "ED7D31", # accent2
    "A5A5A5", # accent3
    "FFC000", # accent4
    "4472C4", # accent5
    "70AD47", # accent6
    "0563C1", # hlink
    "954F72"  # folHlink
)

$slide = $p.Slides.Item(1)
$themeColors = $slide.ThemeColorScheme

for ($i = 0; $i -lt $officeColors.Count; $i++) {
    $themeColors.Item($i + 1).RGB = ConvertTo-ComRGB $officeColors[$i]
}
